# Project Batch 30 Dec 2023
# Adds a new attendance column (D) for 30-Dec-2023, marks each student as
# Present/Reason, documents the two "Reason" entries with review comments,
# widens column D to fit its new "Reason" header, and leaves the selection
# where the author left off (I12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New date header for the 30-Dec-2023 batch (Excel serial 45290) ---
$ws.Range("D1").Value = 45290

# --- Attendance marks for the new date column ---
$ws.Range("D2").Value = "Present"
$ws.Range("D3").Value = "Reason"
$ws.Range("D4").Value = "Present"
$ws.Range("D5").Value = "Present"
$ws.Range("D6").Value = "Present"
$ws.Range("D7").Value = "Reason"
$ws.Range("D8").Value = "Present"

# --- Column D now holds real content ("Reason" header) - size it to fit ---
$ws.Columns("D").ColumnWidth = 10.43

# --- Explain the two "Reason" marks with cell comments ---
$noteD3 = "Hp:" + [char]10 + "Family gate together & Going to Temple"
$commentD3 = $ws.Range("D3").AddComment($noteD3)
$commentD3.Author = "Hp"

$noteD7 = "Hp:" + [char]10 + "Go to the Affidavits Medical Certificate & Non Creamy layer Certificate" + [char]10 + "She not attend Lecture"
$commentD7 = $ws.Range("D7").AddComment($noteD7)
$commentD7.Author = "Hp"

# --- Leave the selection where the author left it ---
[void]$ws.Range("I12").Select()
